# Preprocessing audit: reran updated LASSO model after fix.
# Update the LASSO row (row 2) results: Youden's J, AUC, Sensitivity, Specificity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.4
$ws.Range("E2").Value = "0.79 [0.57–0.95]"
$ws.Range("F2").Value = "66% [33%–100%]"
$ws.Range("G2").Value = "74% [62%–82%]"
